$d = $word.ActiveDocument

# The target paragraph is the final (empty) paragraph in the document body.
$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">A problem regenerative braking Is </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>that current is limited to safe operating region of the batteries, as high currents during regenerative braking decreases the life of the batteries [8]-[9].</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> A solution to this issue is using flywheel and ultracapcitors for regenerative braking as discussed in this paper by S.Bhurse and A.Bhole which conclues that this combination would lead to an increase of range by 16.25%</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>[x1]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>There are man</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">y regenerative braking strategies, several papers are published comparing different approaches[x2]-[x3]. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>The paper by W.Zhang, J.Yang, W.Zhang and Ma[x4]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> compares four different regenerative braking control strategies</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> for Pure Electric Mining Dump Truck.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>Vehicle Speed based control strategy</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> which is parallel braking strategy </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">in </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">which the regenerative braking force </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">increases </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>as the speed of the vehicle increases. This seems a bad approach for</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> braking due to the fact that at high speeds, to maintain braking controllability for safe operation mechanical braking should be engaged. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">The paper restricts the speed of the vehicle to 15km/hr. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">I-Curve based control strategy this control strategy is maximum driver feel braking strategy, It follows the I-Curve for front and rear braking forces. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">β Line control strategy </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:noProof/>
        </w:rPr>
        <w:t>is based on the distribution of axle braking forces are a constant β ratio.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>The last control strategy is F</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>f</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>max</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve"> based strategy. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>The paper concludes that F</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>fmax</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>based strategy improves braking energy recovery compared to other strategies.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>References:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:br/>
        <w:t xml:space="preserve">[8] </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>Pay, S., Baghzouz, Y., “Effectiveness of battery-supercapacitor combination in electric vehicles,” Power Tech Conference Proceedings, 2003 IEEE Bologna, vol.3, no., pp. 6 pp. Vol.3, 23-26 June 2003</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">[9] </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>Gagliardi, F., Pagano, M., “Experimental results of on-board battery</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>ultracapacitor system for electric vehicle applications,” Industrial Electronics, 2002. ISIE 2002. Proceedings of the 2002 IEEE International Symposium on, vol.1, no., pp. 93- 98 vol.1, 2002</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t xml:space="preserve">[x1] </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <w:t>A Review of Regenerative Braking in Electric Vehicles Sneha S. Bhurse A.A. Bhole</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:r>
        <w:t>[x</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>4]</w:t>
      </w:r>
      <w:r>
        <w:t>Research</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> on Regenerative Braking of Pure Electric Mining Dump Truck Wei Zhang * , </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Jue</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Yang , </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Wenming</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Zhang and Fei Ma</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:r>
        <w:t xml:space="preserve">[x2] </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Switched Robust Control of Regenerative Braking of Electric Vehicles </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Xie</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Jing1, 2, Cao Binggang</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>1 ,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Zhang Huarong1 , Xu Dan</w:t>
      </w:r>
    </w:p>
'@

$insertPoint.InsertXML($newContentXml)

# Now re-fetch the (now shifted) final paragraph and strip its stored run
# formatting back down to just noProof, matching the target paragraph mark.
$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$markRange = $targetPara.Range

$markXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr></w:p>'
$markRange.InsertXML($markXml)

Write-Host "Paragraph count:" $d.Paragraphs.Count
